# fix: ajuste de saldo
#
# - Uppercase the two section titles ("Receitas" -> "RECEITAS",
#   "Despesas" -> "DESPESAS").
# - Clarify the "Valor" header as "Valor (BRL)" (used by both the
#   Receitas and Despesas tables).
# - Apply a Brazilian-real currency number format to the "Valor" column.
# - Widen column D to fit the new currency format.
# - Correct a handful of revenue/expense values.
# - Fix the Saldo formula to reference the Total cells in column D
#   instead of column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section headings -------------------------------------------------
$ws.Range("A4").Value = "RECEITAS"
$ws.Range("A12").Value = "DESPESAS"

# --- Column header text -------------------------------------------------
$ws.Range("D5").Value = "Valor (BRL)"
$ws.Range("D13").Value = "Valor (BRL)"

# --- Corrected values ---------------------------------------------------
$ws.Range("D8").Value = 12500
$ws.Range("D10").Value = 6000
$ws.Range("D16").Value = 10000
$ws.Range("D17").Value = 8000

# --- Saldo formula: use the Total cells in column D, not column B -------
$ws.Range("D22").Formula = "=D20 - D21"

# --- Currency number format (R$ #,##0.00) for the whole "Valor" column --
$fmt = "R$ #,##0.00"
$ws.Range("D1").NumberFormat = $fmt
$ws.Range("D5").NumberFormat = $fmt
$ws.Range("D6:D10").NumberFormat = $fmt
$ws.Range("D13").NumberFormat = $fmt
$ws.Range("D14:D17").NumberFormat = $fmt
$ws.Range("D19").NumberFormat = $fmt
$ws.Range("D20:D22").NumberFormat = $fmt

# --- Widen column D so the currency values fit --------------------------
$ws.Columns.Item(4).ColumnWidth = 11.357171694437662
